# NMCARS-PART-5219.docx edit script
# - Adds/changes w:pStyle (List2 / List3) on a number of paragraphs that were
#   plain "Normal" (or had ad-hoc direct pPr formatting) in the "before" doc.
# - Splits a handful of runs that start with a "(i)"/"(ii)"/"(2)" marker into
#   two separate runs: one holding just the marker, one holding the rest of
#   the sentence (mirrors how the author's edit shows up in the OOXML diff).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "(10)(B) In accordance with ASN(RDA) memorandum ..." -> List2
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(10)(B) In accordance with ASN(RDA)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs.Item(1).Style = "List2"

# ---------------------------------------------------------------------
# 2) "(10)(A) The review requirements are not applicable to:" -> List2
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(10)(A) The review requirements are not applicable to:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs.Item(1).Style = "List2"

# ---------------------------------------------------------------------
# 3) "(i) orders placed against single award ..." -> split run after "(i)"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(i) orders placed against single award", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng.Start + 3
$sub = $d.Range($rng.Start, $splitPos)
$sub.Font.Bold = 1
$sub.Font.Bold = 0

# ---------------------------------------------------------------------
# 4) "(ii) awards to small business concerns under the SBIR Program." ->
#    split run after "(ii)"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(ii) awards to small business concerns", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng.Start + 4
$sub = $d.Range($rng.Start, $splitPos)
$sub.Font.Bold = 1
$sub.Font.Bold = 0

# ---------------------------------------------------------------------
# 5) "(1) brief the appointing authority quarterly ..." -> List2
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(1) brief the appointing authority quarterly", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs.Item(1).Style = "List2"

# ---------------------------------------------------------------------
# 6) "(2) conduct SBP program training sessions ..." -> List2
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(2) conduct SBP program training sessions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs.Item(1).Style = "List2"

# ---------------------------------------------------------------------
# 7) "(2)(i) Contracting officers, in evaluating SBA requests ..." ->
#    List1 -> List2, and split run into "(2)" + "(i) Contracting officers..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(2)(i) Contracting officers, in evaluating SBA requests", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs.Item(1).Style = "List2"
$splitPos = $rng.Start + 3
$sub = $d.Range($rng.Start, $splitPos)
$sub.Font.Bold = 1
$sub.Font.Bold = 0

# ---------------------------------------------------------------------
# 8) "(ii) Where possible, procurement activity should be suspended ..." ->
#    List2 -> List3, and split run into "(ii)" + " Where possible..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(ii) Where possible, procurement activity", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Paragraphs.Item(1).Style = "List3"
$splitPos = $rng.Start + 4
$sub = $d.Range($rng.Start, $splitPos)
$sub.Font.Bold = 1
$sub.Font.Bold = 0

# ---------------------------------------------------------------------
# 9) The six "lack of ..." / "deficiencies ..." paragraphs: replace the
#    direct widowControl/tabs pPr formatting with pStyle=List2.
# ---------------------------------------------------------------------
$phrases = @(
    "(1) lack of knowledge and understanding of the work to be performed;",
    "(2) lack of experience in performing requirements of similar size and scope;",
    "(3) lack of resources that are available",
    "(4) lack of ability to comply with subcontracting limitation provisions",
    "(5) lack of ability to meet delivery schedules; and/or",
    "(6) deficiencies in record of performance."
)
foreach ($phrase in $phrases) {
    $rng = $d.Content
    $rng.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Paragraphs.Item(1).Style = "List2"
}
